# Normalize the "Recorded By" column (column G) so that whenever the
# value "System" appears among the comma-separated list of recorders,
# it is moved to the front of the list while the remaining entries
# keep their original relative order. Cells that do not contain
# "System" are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # Column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = @($val -split ", ")

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newVal = $newParts -join ", "

        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
